# Auto-generated edit script applying the diff changes to Alpha_Profits workbook
$wb = $excel.ActiveWorkbook

# ---- Sheet 1: ALC ----
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(32, 8).Value = 16673748
$ws.Cells.Item(32, 10).Value = 25007042
$ws.Cells.Item(32, 12).Value = 25007042
$ws.Cells.Item(32, 14).Value = -25007694
$ws.Cells.Item(70, 8).Value = 3624.5
$ws.Cells.Item(70, 9).Value = 1929.6666
$ws.Cells.Item(70, 10).Value = 4189.4443
$ws.Cells.Item(70, 11).Value = 5788.9998
$ws.Cells.Item(70, 12).Value = 12568.3329
$ws.Cells.Item(70, 13).Value = -5518.9998
$ws.Cells.Item(70, 14).Value = -13108.3329
$ws.Cells.Item(73, 8).Value = 3624.5
$ws.Cells.Item(73, 9).Value = 1929.6666
$ws.Cells.Item(73, 10).Value = 4189.4443
$ws.Cells.Item(73, 11).Value = 5788.9998
$ws.Cells.Item(73, 12).Value = 12568.3329
$ws.Cells.Item(73, 13).Value = -4852.9998
$ws.Cells.Item(73, 14).Value = -14440.3329
$ws.Cells.Item(76, 8).Value = 5142
$ws.Cells.Item(76, 10).Value = 5701
$ws.Cells.Item(76, 12).Value = 5701
$ws.Cells.Item(76, 14).Value = -6331
$ws.Cells.Item(79, 8).Value = 5142
$ws.Cells.Item(79, 10).Value = 5701
$ws.Cells.Item(79, 12).Value = 5701
$ws.Cells.Item(79, 14).Value = -7885
$ws.Cells.Item(97, 8).Value = 1222
$ws.Cells.Item(97, 10).Value = 1222
$ws.Cells.Item(97, 12).Value = 3666
$ws.Cells.Item(97, 14).Value = -4658
$ws.Cells.Item(101, 8).Value = 22076
$ws.Cells.Item(101, 9).Value = 1283.1666
$ws.Cells.Item(101, 10).Value = 63661.668
$ws.Cells.Item(101, 11).Value = 3849.4998
$ws.Cells.Item(101, 12).Value = 190985.004
$ws.Cells.Item(101, 13).Value = -2227.4998
$ws.Cells.Item(101, 14).Value = -194229.004
$ws.Cells.Item(125, 8).Value = 1025
$ws.Cells.Item(125, 9).Value = 777
$ws.Cells.Item(125, 10).Value = 1087
$ws.Cells.Item(125, 11).Value = 6993
$ws.Cells.Item(125, 12).Value = 9783
$ws.Cells.Item(125, 13).Value = -4533
$ws.Cells.Item(125, 14).Value = -14703
$ws.Cells.Item(132, 8).Value = 847.5161000000001
$ws.Cells.Item(132, 9).Value = 825.76666
$ws.Cells.Item(132, 11).Value = 2477.29998
$ws.Cells.Item(132, 13).Value = 52.70002000000022
$ws.Cells.Item(141, 8).Value = 51233.05
$ws.Cells.Item(141, 9).Value = 59633.176
$ws.Cells.Item(141, 10).Value = 3632.3333
$ws.Cells.Item(141, 11).Value = 178899.528
$ws.Cells.Item(141, 12).Value = 10896.9999
$ws.Cells.Item(141, 13).Value = -173719.528
$ws.Cells.Item(141, 14).Value = -21256.9999

# ---- Sheet 2: ARM ----
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(4, 8).Value = 282.16666
$ws.Cells.Item(4, 9).Value = 309.6
$ws.Cells.Item(4, 10).Value = 145
$ws.Cells.Item(4, 11).Value = 309.6
$ws.Cells.Item(4, 12).Value = 145
$ws.Cells.Item(4, 13).Value = -193.6
$ws.Cells.Item(4, 14).Value = -377
$ws.Cells.Item(5, 8).Value = 237.64706
$ws.Cells.Item(5, 10).Value = 149.2
$ws.Cells.Item(5, 12).Value = 149.2
$ws.Cells.Item(5, 14).Value = -373.2
$ws.Cells.Item(24, 8).Value = 31000
$ws.Cells.Item(24, 10).Value = 31000
$ws.Cells.Item(24, 12).Value = 31000
$ws.Cells.Item(24, 14).Value = -31748
$ws.Cells.Item(26, 8).Value = 4323.75
$ws.Cells.Item(26, 9).Value = 4333.3335
$ws.Cells.Item(26, 11).Value = 4333.3335
$ws.Cells.Item(26, 13).Value = -4003.3335
$ws.Cells.Item(32, 8).Value = 5568.857
$ws.Cells.Item(32, 9).Value = 5568.857
$ws.Cells.Item(32, 11).Value = 5568.857
$ws.Cells.Item(32, 13).Value = -5281.857
$ws.Cells.Item(33, 8).Value = 8500
$ws.Cells.Item(33, 9).Value = 8500
$ws.Cells.Item(33, 11).Value = 8500
$ws.Cells.Item(33, 13).Value = -8171
$ws.Cells.Item(43, 8).Value = 22175.334
$ws.Cells.Item(43, 9).Value = 14930
$ws.Cells.Item(43, 10).Value = 23624.4
$ws.Cells.Item(43, 11).Value = 14930
$ws.Cells.Item(43, 12).Value = 23624.4
$ws.Cells.Item(43, 13).Value = -14617
$ws.Cells.Item(43, 14).Value = -24250.4
$ws.Cells.Item(45, 8).Value = 3245.9092
$ws.Cells.Item(45, 10).Value = 5424.6665
$ws.Cells.Item(45, 12).Value = 5424.6665
$ws.Cells.Item(45, 14).Value = -6178.6665
$ws.Cells.Item(61, 8).Value = 4693.5
$ws.Cells.Item(61, 9).Value = 4528.4707
$ws.Cells.Item(61, 11).Value = 4528.4707
$ws.Cells.Item(61, 13).Value = -4316.4707
$ws.Cells.Item(74, 8).Value = 1179.4634
$ws.Cells.Item(74, 9).Value = 1035.3448
$ws.Cells.Item(74, 10).Value = 1527.75
$ws.Cells.Item(74, 11).Value = 1035.3448
$ws.Cells.Item(74, 12).Value = 1527.75
$ws.Cells.Item(74, 13).Value = -161.3448000000001
$ws.Cells.Item(74, 14).Value = -3275.75
$ws.Cells.Item(77, 8).Value = 1179.4634
$ws.Cells.Item(77, 9).Value = 1035.3448
$ws.Cells.Item(77, 10).Value = 1527.75
$ws.Cells.Item(77, 11).Value = 5176.724
$ws.Cells.Item(77, 12).Value = 7638.75
$ws.Cells.Item(77, 13).Value = -808.7240000000002
$ws.Cells.Item(77, 14).Value = -16374.75
$ws.Cells.Item(100, 8).Value = 31000
$ws.Cells.Item(100, 10).Value = 31000
$ws.Cells.Item(100, 12).Value = 31000
$ws.Cells.Item(100, 14).Value = -33164
$ws.Cells.Item(105, 8).Value = 0
$ws.Cells.Item(105, 10).Value = 0
$ws.Cells.Item(105, 12).Value = 0
$ws.Cells.Item(105, 14).ClearContents()
$ws.Cells.Item(132, 8).Value = 2945
$ws.Cells.Item(132, 9).Value = 2945
$ws.Cells.Item(132, 11).Value = 8835
$ws.Cells.Item(132, 13).Value = -6305
$ws.Cells.Item(136, 8).Value = 4693.5
$ws.Cells.Item(136, 9).Value = 4528.4707
$ws.Cells.Item(136, 11).Value = 13585.4121
$ws.Cells.Item(136, 13).Value = -11035.4121

# ---- Sheet 3: BSM ----
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(4, 8).Value = 237.64706
$ws.Cells.Item(4, 10).Value = 149.2
$ws.Cells.Item(4, 12).Value = 149.2
$ws.Cells.Item(4, 14).Value = -379.2
$ws.Cells.Item(13, 8).Value = 86500
$ws.Cells.Item(13, 9).Value = 23000
$ws.Cells.Item(13, 11).Value = 23000
$ws.Cells.Item(13, 13).Value = -22832
$ws.Cells.Item(25, 8).Value = 8715.857
$ws.Cells.Item(25, 9).Value = 3602.2
$ws.Cells.Item(25, 10).Value = 21500
$ws.Cells.Item(25, 11).Value = 3602.2
$ws.Cells.Item(25, 12).Value = 21500
$ws.Cells.Item(25, 13).Value = -3367.2
$ws.Cells.Item(25, 14).Value = -21970
$ws.Cells.Item(105, 8).Value = 2370
$ws.Cells.Item(105, 9).Value = 1918
$ws.Cells.Item(105, 11).Value = 1918
$ws.Cells.Item(105, 13).Value = -171
$ws.Cells.Item(134, 8).Value = 3374.4333
$ws.Cells.Item(134, 9).Value = 3463.2068
$ws.Cells.Item(134, 11).Value = 10389.6204
$ws.Cells.Item(134, 13).Value = -7854.6204

# ---- Sheet 4: CRP ----
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(7, 8).Value = 2976300.8
$ws.Cells.Item(7, 9).Value = 3472329.2
$ws.Cells.Item(7, 11).Value = 3472329.2
$ws.Cells.Item(7, 13).Value = -3472216.2
$ws.Cells.Item(23, 8).Value = 16405.857
$ws.Cells.Item(23, 9).Value = 8498.75
$ws.Cells.Item(23, 11).Value = 8498.75
$ws.Cells.Item(23, 13).Value = -8258.75
$ws.Cells.Item(27, 8).Value = 16405.857
$ws.Cells.Item(27, 9).Value = 8498.75
$ws.Cells.Item(27, 11).Value = 8498.75
$ws.Cells.Item(27, 13).Value = -8306.75
$ws.Cells.Item(31, 8).Value = 1798.909
$ws.Cells.Item(31, 9).Value = 1936.1111
$ws.Cells.Item(31, 10).Value = 1181.5
$ws.Cells.Item(31, 11).Value = 1936.1111
$ws.Cells.Item(31, 12).Value = 1181.5
$ws.Cells.Item(31, 13).Value = -1641.1111
$ws.Cells.Item(31, 14).Value = -1771.5
$ws.Cells.Item(32, 8).Value = 510
$ws.Cells.Item(32, 9).Value = 510
$ws.Cells.Item(32, 11).Value = 510
$ws.Cells.Item(32, 13).Value = -194
$ws.Cells.Item(34, 8).Value = 1798.909
$ws.Cells.Item(34, 9).Value = 1936.1111
$ws.Cells.Item(34, 10).Value = 1181.5
$ws.Cells.Item(34, 11).Value = 1936.1111
$ws.Cells.Item(34, 12).Value = 1181.5
$ws.Cells.Item(34, 13).Value = -1734.1111
$ws.Cells.Item(34, 14).Value = -1585.5
$ws.Cells.Item(36, 8).Value = 17012
$ws.Cells.Item(36, 9).Value = 10024
$ws.Cells.Item(36, 11).Value = 10024
$ws.Cells.Item(36, 13).Value = -9636
$ws.Cells.Item(40, 8).Value = 17012
$ws.Cells.Item(40, 9).Value = 10024
$ws.Cells.Item(40, 11).Value = 10024
$ws.Cells.Item(40, 13).Value = -9864
$ws.Cells.Item(132, 8).Value = 1282.4286
$ws.Cells.Item(132, 9).Value = 997
$ws.Cells.Item(132, 11).Value = 2991
$ws.Cells.Item(132, 13).Value = -461

# ---- Sheet 5: CUL ----
$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(4, 8).Value = 376791.12
$ws.Cells.Item(4, 9).Value = 443.33334
$ws.Cells.Item(4, 10).Value = 602599.8
$ws.Cells.Item(4, 11).Value = 1330.00002
$ws.Cells.Item(4, 12).Value = 1807799.4
$ws.Cells.Item(4, 13).Value = -1218.00002
$ws.Cells.Item(4, 14).Value = -1808023.4
$ws.Cells.Item(39, 8).Value = 5163.636
$ws.Cells.Item(39, 10).Value = 5670
$ws.Cells.Item(39, 12).Value = 17010
$ws.Cells.Item(39, 14).Value = -17598
$ws.Cells.Item(107, 8).Value = 520.3333
$ws.Cells.Item(107, 10).Value = 507.30768
$ws.Cells.Item(107, 12).Value = 1521.92304
$ws.Cells.Item(107, 14).Value = -5361.92304
$ws.Cells.Item(117, 8).Value = 1932.25
$ws.Cells.Item(117, 9).Value = 1932.25
$ws.Cells.Item(117, 11).Value = 5796.75
$ws.Cells.Item(117, 13).Value = -2354.75
$ws.Cells.Item(132, 8).Value = 10969.857
$ws.Cells.Item(132, 10).Value = 14217.8
$ws.Cells.Item(132, 12).Value = 127960.2
$ws.Cells.Item(132, 14).Value = -133020.2

# ---- Sheet 6: GSM ----
$ws = $wb.Worksheets.Item(6)
$ws.Cells.Item(18, 8).Value = 0
$ws.Cells.Item(18, 10).Value = 0
$ws.Cells.Item(18, 12).Value = 0
$ws.Cells.Item(18, 14).ClearContents()
$ws.Cells.Item(20, 8).Value = 10000000
$ws.Cells.Item(20, 10).Value = 0
$ws.Cells.Item(20, 12).Value = 0
$ws.Cells.Item(20, 14).ClearContents()
$ws.Cells.Item(24, 8).Value = 5015500
$ws.Cells.Item(24, 10).Value = 31000
$ws.Cells.Item(24, 12).Value = 31000
$ws.Cells.Item(24, 14).Value = -31346
$ws.Cells.Item(29, 8).Value = 12330
$ws.Cells.Item(29, 9).Value = 6990
$ws.Cells.Item(29, 11).Value = 6990
$ws.Cells.Item(29, 13).Value = -6700
$ws.Cells.Item(69, 8).Value = 61636.184
$ws.Cells.Item(69, 9).Value = 0
$ws.Cells.Item(69, 10).Value = 61636.184
$ws.Cells.Item(69, 11).Value = 0
$ws.Cells.Item(69, 12).Value = 61636.184
$ws.Cells.Item(69, 13).ClearContents()
$ws.Cells.Item(69, 14).Value = -63134.184
$ws.Cells.Item(72, 8).Value = 61636.184
$ws.Cells.Item(72, 9).Value = 0
$ws.Cells.Item(72, 10).Value = 61636.184
$ws.Cells.Item(72, 11).Value = 0
$ws.Cells.Item(72, 12).Value = 184908.552
$ws.Cells.Item(72, 13).ClearContents()
$ws.Cells.Item(72, 14).Value = -192396.552
$ws.Cells.Item(80, 8).Value = 8279.166999999999
$ws.Cells.Item(80, 9).Value = 1945
$ws.Cells.Item(80, 11).Value = 1945
$ws.Cells.Item(80, 13).Value = -947
$ws.Cells.Item(83, 8).Value = 8279.166999999999
$ws.Cells.Item(83, 9).Value = 1945
$ws.Cells.Item(83, 11).Value = 9725
$ws.Cells.Item(83, 13).Value = -4733
$ws.Cells.Item(102, 8).Value = 947.8421
$ws.Cells.Item(102, 9).Value = 875.2143
$ws.Cells.Item(102, 11).Value = 875.2143
$ws.Cells.Item(102, 13).Value = 746.7857
$ws.Cells.Item(132, 8).Value = 3291.5
$ws.Cells.Item(132, 9).Value = 3816.3333
$ws.Cells.Item(132, 10).Value = 2766.6667
$ws.Cells.Item(132, 11).Value = 11448.9999
$ws.Cells.Item(132, 12).Value = 8300.000100000001
$ws.Cells.Item(132, 13).Value = -8918.999899999999
$ws.Cells.Item(132, 14).Value = -13360.0001

# ---- Sheet 7: LTW ----
$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(7, 8).Value = 3757.7273
$ws.Cells.Item(7, 9).Value = 3491
$ws.Cells.Item(7, 11).Value = 3491
$ws.Cells.Item(7, 13).Value = -3379
$ws.Cells.Item(29, 8).Value = 0
$ws.Cells.Item(29, 9).Value = 0
$ws.Cells.Item(29, 11).Value = 0
$ws.Cells.Item(29, 13).ClearContents()
$ws.Cells.Item(31, 8).Value = 1988.3334
$ws.Cells.Item(31, 9).Value = 619
$ws.Cells.Item(31, 10).Value = 14997
$ws.Cells.Item(31, 11).Value = 619
$ws.Cells.Item(31, 12).Value = 14997
$ws.Cells.Item(31, 13).Value = -371
$ws.Cells.Item(31, 14).Value = -15493
$ws.Cells.Item(40, 8).Value = 6107.1763
$ws.Cells.Item(40, 9).Value = 4762
$ws.Cells.Item(40, 11).Value = 4762
$ws.Cells.Item(40, 13).Value = -4626
$ws.Cells.Item(76, 8).Value = 0
$ws.Cells.Item(76, 10).Value = 0
$ws.Cells.Item(76, 12).Value = 0
$ws.Cells.Item(76, 14).ClearContents()
$ws.Cells.Item(79, 8).Value = 0
$ws.Cells.Item(79, 10).Value = 0
$ws.Cells.Item(79, 12).Value = 0
$ws.Cells.Item(79, 14).ClearContents()
$ws.Cells.Item(126, 8).Value = 3757.7273
$ws.Cells.Item(126, 9).Value = 3491
$ws.Cells.Item(126, 11).Value = 10473
$ws.Cells.Item(126, 13).Value = -8003
$ws.Cells.Item(132, 8).Value = 3799.6667
$ws.Cells.Item(132, 9).Value = 3400
$ws.Cells.Item(132, 10).Value = 3999.5
$ws.Cells.Item(132, 11).Value = 10200
$ws.Cells.Item(132, 12).Value = 11998.5
$ws.Cells.Item(132, 13).Value = -7670
$ws.Cells.Item(132, 14).Value = -17058.5
$ws.Cells.Item(136, 8).Value = 5166.5654
$ws.Cells.Item(136, 9).Value = 3719.0557
$ws.Cells.Item(136, 10).Value = 10377.6
$ws.Cells.Item(136, 11).Value = 11157.1671
$ws.Cells.Item(136, 12).Value = 31132.8
$ws.Cells.Item(136, 13).Value = -8607.167099999999
$ws.Cells.Item(136, 14).Value = -36232.8

# ---- Sheet 8: WVR ----
$ws = $wb.Worksheets.Item(8)
$ws.Cells.Item(26, 8).Value = 0
$ws.Cells.Item(26, 9).Value = 0
$ws.Cells.Item(26, 11).Value = 0
$ws.Cells.Item(26, 13).ClearContents()
$ws.Cells.Item(29, 8).Value = 7182
$ws.Cells.Item(29, 10).Value = 4000
$ws.Cells.Item(29, 12).Value = 4000
$ws.Cells.Item(29, 14).Value = -4580
$ws.Cells.Item(31, 8).Value = 30000
$ws.Cells.Item(31, 10).Value = 30000
$ws.Cells.Item(31, 12).Value = 30000
$ws.Cells.Item(31, 14).Value = -30696
$ws.Cells.Item(32, 8).Value = 11665.5
$ws.Cells.Item(32, 9).Value = 3331
$ws.Cells.Item(32, 11).Value = 3331
$ws.Cells.Item(32, 13).Value = -3014
$ws.Cells.Item(34, 8).Value = 23255.334
$ws.Cells.Item(34, 9).Value = 24439
$ws.Cells.Item(34, 10).Value = 20888
$ws.Cells.Item(34, 11).Value = 24439
$ws.Cells.Item(34, 12).Value = 20888
$ws.Cells.Item(34, 13).Value = -24236
$ws.Cells.Item(34, 14).Value = -21294
$ws.Cells.Item(82, 8).Value = 37777.5
$ws.Cells.Item(82, 10).Value = 37777.5
$ws.Cells.Item(82, 12).Value = 37777.5
$ws.Cells.Item(82, 14).Value = -38543.5
$ws.Cells.Item(85, 8).Value = 37777.5
$ws.Cells.Item(85, 10).Value = 37777.5
$ws.Cells.Item(85, 12).Value = 37777.5
$ws.Cells.Item(85, 14).Value = -40429.5
$ws.Cells.Item(111, 8).Value = 50000
$ws.Cells.Item(111, 10).Value = 50000
$ws.Cells.Item(111, 12).Value = 50000
$ws.Cells.Item(111, 14).Value = -58180
$ws.Cells.Item(113, 8).Value = 1683.8948
$ws.Cells.Item(113, 10).Value = 2258.0833
$ws.Cells.Item(113, 12).Value = 6774.249899999999
$ws.Cells.Item(113, 14).Value = -11114.2499
$ws.Cells.Item(122, 8).Value = 2385.963
$ws.Cells.Item(122, 9).Value = 2326.3333
$ws.Cells.Item(122, 10).Value = 2594.6667
$ws.Cells.Item(122, 11).Value = 6978.999899999999
$ws.Cells.Item(122, 12).Value = 7784.000100000001
$ws.Cells.Item(122, 13).Value = -4528.999899999999
$ws.Cells.Item(122, 14).Value = -12684.0001
$ws.Cells.Item(126, 8).Value = 2540.2
$ws.Cells.Item(126, 9).Value = 1600.75
$ws.Cells.Item(126, 11).Value = 4802.25
$ws.Cells.Item(126, 13).Value = -2332.25
$ws.Cells.Item(132, 8).Value = 3010.634
$ws.Cells.Item(132, 9).Value = 2480.3845
$ws.Cells.Item(132, 10).Value = 3929.7334
$ws.Cells.Item(132, 11).Value = 7441.1535
$ws.Cells.Item(132, 12).Value = 11789.2002
$ws.Cells.Item(132, 13).Value = -4911.1535
$ws.Cells.Item(132, 14).Value = -16849.2002
$ws.Cells.Item(136, 8).Value = 690.69696
$ws.Cells.Item(136, 9).Value = 586.1667
$ws.Cells.Item(136, 11).Value = 1758.5001
$ws.Cells.Item(136, 13).Value = 791.4999
